$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(19, 8).Value = 2071.2856  # H19: 2071.4285 -> 2071.2856
$ws.Cells.Item(19, 9).Value = 1749.75  # I19: 1750 -> 1749.75
$ws.Cells.Item(19, 11).Value = 1749.75  # K19: 1750 -> 1749.75
$ws.Cells.Item(19, 13).Value = -1574.75  # M19: -1575 -> -1574.75
$ws.Cells.Item(33, 8).Value = 725.375  # H33: 704.7692 -> 725.375
$ws.Cells.Item(33, 9).Value = 467.16666  # I33: 366.2 -> 467.16666
$ws.Cells.Item(33, 10).Value = 1500  # J33: 1833.3334 -> 1500
$ws.Cells.Item(33, 11).Value = 467.16666  # K33: 366.2 -> 467.16666
$ws.Cells.Item(33, 12).Value = 1500  # L33: 1833.3334 -> 1500
$ws.Cells.Item(33, 13).Value = -238.16666  # M33: -137.2 -> -238.16666
$ws.Cells.Item(33, 14).Value = -1958  # N33: -2291.3334 -> -1958
$ws.Cells.Item(106, 8).Value = 1365  # H106: 1331.85 -> 1365
$ws.Cells.Item(106, 9).Value = 1478  # I106: 1418.3077 -> 1478
$ws.Cells.Item(106, 11).Value = 1478  # K106: 1418.3077 -> 1478
$ws.Cells.Item(106, 13).Value = -847  # M106: -787.3077000000001 -> -847
$ws.Cells.Item(113, 8).Value = 7446.7896  # H113: 7382.722 -> 7446.7896
$ws.Cells.Item(113, 9).Value = 6688  # I113: 6414.857 -> 6688
$ws.Cells.Item(113, 11).Value = 6688  # K113: 6414.857 -> 6688
$ws.Cells.Item(113, 13).Value = -3434  # M113: -3160.857 -> -3434
$ws.Cells.Item(125, 8).Value = 2939.5  # H125: 2639.5557 -> 2939.5
$ws.Cells.Item(125, 9).Value = 769.5  # I125: 693.8570999999999 -> 769.5
$ws.Cells.Item(125, 11).Value = 6925.5  # K125: 6244.7139 -> 6925.5
$ws.Cells.Item(125, 13).Value = -4465.5  # M125: -3784.7139 -> -4465.5
$ws.Cells.Item(132, 8).Value = 10597.693  # H132: 10491.381 -> 10597.693
$ws.Cells.Item(132, 9).Value = 1750.9231  # I132: 1791.4717 -> 1750.9231
$ws.Cells.Item(132, 11).Value = 5252.7693  # K132: 5374.4151 -> 5252.7693
$ws.Cells.Item(132, 13).Value = -2722.7693  # M132: -2844.4151 -> -2722.7693
$ws.Cells.Item(137, 8).Value = 2510.15  # H137: 2513.5122 -> 2510.15
$ws.Cells.Item(137, 9).Value = 2322.054  # I137: 2330.6316 -> 2322.054
$ws.Cells.Item(137, 11).Value = 6966.162  # K137: 6991.8948 -> 6966.162
$ws.Cells.Item(137, 13).Value = -4416.162  # M137: -4441.8948 -> -4416.162

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 14103.04  # H2: 14728.167 -> 14103.04
$ws.Cells.Item(2, 9).Value = 16475.523  # I2: 17284.3 -> 16475.523
$ws.Cells.Item(2, 10).Value = 1647.5  # J2: 1947.5 -> 1647.5
$ws.Cells.Item(2, 11).Value = 16475.523  # K2: 17284.3 -> 16475.523
$ws.Cells.Item(2, 12).Value = 1647.5  # L2: 1947.5 -> 1647.5
$ws.Cells.Item(2, 13).Value = -16362.523  # M2: -17171.3 -> -16362.523
$ws.Cells.Item(2, 14).Value = -1873.5  # N2: -2173.5 -> -1873.5
$ws.Cells.Item(32, 8).Value = 11115991  # H32: 11116010 -> 11115991
$ws.Cells.Item(32, 9).Value = 11368491  # I32: 11368510 -> 11368491
$ws.Cells.Item(32, 11).Value = 11368491  # K32: 11368510 -> 11368491
$ws.Cells.Item(32, 13).Value = -11368204  # M32: -11368223 -> -11368204
$ws.Cells.Item(116, 8).Value = 14103.04  # H116: 14728.167 -> 14103.04
$ws.Cells.Item(116, 9).Value = 16475.523  # I116: 17284.3 -> 16475.523
$ws.Cells.Item(116, 10).Value = 1647.5  # J116: 1947.5 -> 1647.5
$ws.Cells.Item(116, 11).Value = 16475.523  # K116: 17284.3 -> 16475.523
$ws.Cells.Item(116, 12).Value = 1647.5  # L116: 1947.5 -> 1647.5
$ws.Cells.Item(116, 13).Value = -14181.523  # M116: -14990.3 -> -14181.523
$ws.Cells.Item(116, 14).Value = -6235.5  # N116: -6535.5 -> -6235.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 14103.04  # H3: 14728.167 -> 14103.04
$ws.Cells.Item(3, 9).Value = 16475.523  # I3: 17284.3 -> 16475.523
$ws.Cells.Item(3, 10).Value = 1647.5  # J3: 1947.5 -> 1647.5
$ws.Cells.Item(3, 11).Value = 16475.523  # K3: 17284.3 -> 16475.523
$ws.Cells.Item(3, 12).Value = 1647.5  # L3: 1947.5 -> 1647.5
$ws.Cells.Item(3, 13).Value = -16361.523  # M3: -17170.3 -> -16361.523
$ws.Cells.Item(3, 14).Value = -1875.5  # N3: -2175.5 -> -1875.5
$ws.Cells.Item(64, 8).Value = 912.6667  # H64: 972.44446 -> 912.6667
$ws.Cells.Item(64, 9).Value = 838.75  # I64: 818.3333 -> 838.75
$ws.Cells.Item(64, 10).Value = 949.625  # J64: 1049.5 -> 949.625
$ws.Cells.Item(64, 11).Value = 838.75  # K64: 818.3333 -> 838.75
$ws.Cells.Item(64, 12).Value = 949.625  # L64: 1049.5 -> 949.625
$ws.Cells.Item(64, 13).Value = -613.75  # M64: -593.3333 -> -613.75
$ws.Cells.Item(64, 14).Value = -1399.625  # N64: -1499.5 -> -1399.625
$ws.Cells.Item(67, 8).Value = 912.6667  # H67: 972.44446 -> 912.6667
$ws.Cells.Item(67, 9).Value = 838.75  # I67: 818.3333 -> 838.75
$ws.Cells.Item(67, 10).Value = 949.625  # J67: 1049.5 -> 949.625
$ws.Cells.Item(67, 11).Value = 838.75  # K67: 818.3333 -> 838.75
$ws.Cells.Item(67, 12).Value = 949.625  # L67: 1049.5 -> 949.625
$ws.Cells.Item(67, 13).Value = -58.75  # M67: -38.33330000000001 -> -58.75
$ws.Cells.Item(67, 14).Value = -2509.625  # N67: -2609.5 -> -2509.625
$ws.Cells.Item(86, 8).Value = 2642.7083  # H86: 2661.4167 -> 2642.7083
$ws.Cells.Item(86, 9).Value = 1966.6765  # I86: 1944.7142 -> 1966.6765
$ws.Cells.Item(86, 10).Value = 4284.5  # J86: 4591 -> 4284.5
$ws.Cells.Item(86, 11).Value = 1966.6765  # K86: 1944.7142 -> 1966.6765
$ws.Cells.Item(86, 12).Value = 4284.5  # L86: 4591 -> 4284.5
$ws.Cells.Item(86, 13).Value = -843.6765  # M86: -821.7141999999999 -> -843.6765
$ws.Cells.Item(86, 14).Value = -6530.5  # N86: -6837 -> -6530.5
$ws.Cells.Item(89, 8).Value = 2642.7083  # H89: 2661.4167 -> 2642.7083
$ws.Cells.Item(89, 9).Value = 1966.6765  # I89: 1944.7142 -> 1966.6765
$ws.Cells.Item(89, 10).Value = 4284.5  # J89: 4591 -> 4284.5
$ws.Cells.Item(89, 11).Value = 9833.3825  # K89: 9723.571 -> 9833.3825
$ws.Cells.Item(89, 12).Value = 21422.5  # L89: 22955 -> 21422.5
$ws.Cells.Item(89, 13).Value = -4217.3825  # M89: -4107.571 -> -4217.3825
$ws.Cells.Item(89, 14).Value = -32654.5  # N89: -34187 -> -32654.5
$ws.Cells.Item(94, 9).Value = 2509.2727  # I94: 2408.8696 -> 2509.2727
$ws.Cells.Item(94, 10).Value = 337.5  # J94: 475 -> 337.5
$ws.Cells.Item(94, 11).Value = 2509.2727  # K94: 2408.8696 -> 2509.2727
$ws.Cells.Item(94, 12).Value = 337.5  # L94: 475 -> 337.5
$ws.Cells.Item(94, 13).Value = -2058.2727  # M94: -1957.8696 -> -2058.2727
$ws.Cells.Item(94, 14).Value = -1239.5  # N94: -1377 -> -1239.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2243.0688  # H31: 2243.5862 -> 2243.0688
$ws.Cells.Item(31, 9).Value = 2127.8572  # I31: 2128.5715 -> 2127.8572
$ws.Cells.Item(31, 11).Value = 2127.8572  # K31: 2128.5715 -> 2127.8572
$ws.Cells.Item(31, 13).Value = -1832.8572  # M31: -1833.5715 -> -1832.8572
$ws.Cells.Item(34, 8).Value = 2243.0688  # H34: 2243.5862 -> 2243.0688
$ws.Cells.Item(34, 9).Value = 2127.8572  # I34: 2128.5715 -> 2127.8572
$ws.Cells.Item(34, 11).Value = 2127.8572  # K34: 2128.5715 -> 2127.8572
$ws.Cells.Item(34, 13).Value = -1925.8572  # M34: -1926.5715 -> -1925.8572
$ws.Cells.Item(58, 8).Value = 1943.5  # H58: 1771.5676 -> 1943.5
$ws.Cells.Item(58, 9).Value = 1213.579  # I58: 1140.6522 -> 1213.579
$ws.Cells.Item(58, 10).Value = 3010.3076  # J58: 2808.0715 -> 3010.3076
$ws.Cells.Item(58, 11).Value = 1213.579  # K58: 1140.6522 -> 1213.579
$ws.Cells.Item(58, 12).Value = 3010.3076  # L58: 2808.0715 -> 3010.3076
$ws.Cells.Item(58, 13).Value = -1010.579  # M58: -937.6522 -> -1010.579
$ws.Cells.Item(58, 14).Value = -3416.3076  # N58: -3214.0715 -> -3416.3076
$ws.Cells.Item(132, 8).Value = 2173.7727  # H132: 2239.1428 -> 2173.7727
$ws.Cells.Item(132, 9).Value = 2106.4736  # I132: 2179 -> 2106.4736
$ws.Cells.Item(132, 11).Value = 6319.4208  # K132: 6537 -> 6319.4208
$ws.Cells.Item(132, 13).Value = -3789.4208  # M132: -4007 -> -3789.4208
$ws.Cells.Item(134, 8).Value = 2179.0952  # H134: 2127.5454 -> 2179.0952
$ws.Cells.Item(134, 9).Value = 1769  # I134: 1726.4117 -> 1769
$ws.Cells.Item(134, 11).Value = 5307  # K134: 5179.2351 -> 5307
$ws.Cells.Item(134, 13).Value = -2772  # M134: -2644.2351 -> -2772
$ws.Cells.Item(136, 8).Value = 1943.5  # H136: 1771.5676 -> 1943.5
$ws.Cells.Item(136, 9).Value = 1213.579  # I136: 1140.6522 -> 1213.579
$ws.Cells.Item(136, 10).Value = 3010.3076  # J136: 2808.0715 -> 3010.3076
$ws.Cells.Item(136, 11).Value = 3640.737  # K136: 3421.9566 -> 3640.737
$ws.Cells.Item(136, 12).Value = 9030.9228  # L136: 8424.2145 -> 9030.9228
$ws.Cells.Item(136, 13).Value = -1090.737  # M136: -871.9566 -> -1090.737
$ws.Cells.Item(136, 14).Value = -14130.9228  # N136: -13524.2145 -> -14130.9228

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 583.1667  # H34: 850 -> 583.1667
$ws.Cells.Item(34, 9).Value = 99.8  # I34: 100 -> 99.8
$ws.Cells.Item(34, 10).Value = 3000  # J34: 1600 -> 3000
$ws.Cells.Item(34, 11).Value = 299.4  # K34: 300 -> 299.4
$ws.Cells.Item(34, 12).Value = 9000  # L34: 4800 -> 9000
$ws.Cells.Item(34, 13).Value = -215.4  # M34: -216 -> -215.4
$ws.Cells.Item(34, 14).Value = -9168  # N34: -4968 -> -9168
$ws.Cells.Item(35, 8).Value = 354  # H35: 300 -> 354
$ws.Cells.Item(35, 10).Value = 354  # J35: 300 -> 354
$ws.Cells.Item(35, 12).Value = 1062  # L35: 900 -> 1062
$ws.Cells.Item(35, 14).Value = -1638  # N35: -1476 -> -1638
$ws.Cells.Item(59, 8).Value = 162500  # H59: 211484.67 -> 162500
$ws.Cells.Item(59, 9).Value = 162500  # I59: 157502.5 -> 162500
$ws.Cells.Item(59, 10).Value = 0  # J59: 319449 -> 0
$ws.Cells.Item(59, 11).Value = 487500  # K59: 472507.5 -> 487500
$ws.Cells.Item(59, 12).Value = 0  # L59: 958347 -> 0
$ws.Cells.Item(59, 13).ClearContents()  # M59: was -471967.5
$ws.Cells.Item(59, 14).Value = -486960  # N59: -959427 -> -486960
$ws.Cells.Item(124, 8).Value = 1099.75  # H124: 1749.5 -> 1099.75
$ws.Cells.Item(124, 9).Value = 466.33334  # I124: 1332.6666 -> 466.33334
$ws.Cells.Item(124, 11).Value = 1399.00002  # K124: 3997.9998 -> 1399.00002
$ws.Cells.Item(124, 13).Value = 3510.99998  # M124: 912.0001999999999 -> 3510.99998
$ws.Cells.Item(137, 8).Value = 3247.0557  # H137: 3116.2856 -> 3247.0557
$ws.Cells.Item(137, 9).Value = 2857.1428  # I137: 2812.5 -> 2857.1428
$ws.Cells.Item(137, 10).Value = 3495.182  # J137: 3303.2307 -> 3495.182
$ws.Cells.Item(137, 11).Value = 8571.428400000001  # K137: 8437.5 -> 8571.428400000001
$ws.Cells.Item(137, 12).Value = 10485.546  # L137: 9909.6921 -> 10485.546
$ws.Cells.Item(137, 13).Value = -3471.428400000001  # M137: -3337.5 -> -3471.428400000001
$ws.Cells.Item(137, 14).Value = -20685.546  # N137: -20109.6921 -> -20685.546

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 2657.0588  # H126: 2710.8125 -> 2657.0588
$ws.Cells.Item(126, 9).Value = 2698  # I126: 2762.3572 -> 2698
$ws.Cells.Item(126, 11).Value = 8094  # K126: 8287.071599999999 -> 8094
$ws.Cells.Item(126, 13).Value = -5624  # M126: -5817.071599999999 -> -5624
$ws.Cells.Item(132, 8).Value = 2325.158  # H132: 2244.1904 -> 2325.158
$ws.Cells.Item(132, 9).Value = 1878.9333  # I132: 1831.4117 -> 1878.9333
$ws.Cells.Item(132, 11).Value = 5636.7999  # K132: 5494.2351 -> 5636.7999
$ws.Cells.Item(132, 13).Value = -3106.7999  # M132: -2964.2351 -> -3106.7999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 3600.1462  # H132: 3665.625 -> 3600.1462
$ws.Cells.Item(132, 9).Value = 2196.5938  # I132: 2235.8064 -> 2196.5938
$ws.Cells.Item(132, 11).Value = 6589.7814  # K132: 6707.4192 -> 6589.7814
$ws.Cells.Item(132, 13).Value = -4059.7814  # M132: -4177.4192 -> -4059.7814
$ws.Cells.Item(136, 8).Value = 2830.3044  # H136: 2927.6365 -> 2830.3044
$ws.Cells.Item(136, 9).Value = 2421.2856  # I136: 2507.9 -> 2421.2856
$ws.Cells.Item(136, 11).Value = 7263.8568  # K136: 7523.700000000001 -> 7263.8568
$ws.Cells.Item(136, 13).Value = -4713.8568  # M136: -4973.700000000001 -> -4713.8568

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 4729.4443  # H81: 3338.4285 -> 4729.4443
$ws.Cells.Item(81, 10).Value = 9000  # J81: 3896.625 -> 9000
$ws.Cells.Item(81, 12).Value = 18000  # L81: 7793.25 -> 18000
$ws.Cells.Item(81, 14).Value = -20122  # N81: -9915.25 -> -20122
$ws.Cells.Item(84, 8).Value = 4729.4443  # H84: 3338.4285 -> 4729.4443
$ws.Cells.Item(84, 10).Value = 9000  # J84: 3896.625 -> 9000
$ws.Cells.Item(84, 12).Value = 90000  # L84: 38966.25 -> 90000
$ws.Cells.Item(84, 14).Value = -100608  # N84: -49574.25 -> -100608
$ws.Cells.Item(104, 8).Value = 65185  # H104: 0 -> 65185
$ws.Cells.Item(104, 10).Value = 65185  # J104: 0 -> 65185
$ws.Cells.Item(104, 12).Value = 65185  # L104: 0 -> 65185
$ws.Cells.Item(104, 14).Value = -72173  # N104: None -> -72173
$ws.Cells.Item(113, 8).Value = 656.55554  # H113: 710.5599999999999 -> 656.55554
$ws.Cells.Item(113, 9).Value = 382.61905  # I113: 424.8421 -> 382.61905
$ws.Cells.Item(113, 10).Value = 1147.85715  # J113: 1615.3334 -> 1147.85715
$ws.Cells.Item(113, 11).Value = 1147.85715  # K113: 1274.5263 -> 1147.85715
$ws.Cells.Item(113, 13).Value = 1022.14285  # M113: 895.4737 -> 1022.14285
